$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 16 (3rd shape): switch its table style from the custom
#    "Table_0" style to the other built-in table style.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{12A8EF02-4A23-4DD0-9239-35A2F9C71849}")

# ---------------------------------------------------------------------------
# 2) Swap the presentation theme's colour palette: it currently carries the
#    "Integral" accent colours; replace every slot with the plain "Office
#    Theme" palette colours (same dk1/lt1, new dk2/lt2/accent1-6/hlink/
#    folHlink).
# ---------------------------------------------------------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
# order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#        8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# values are OLE RGB (0x00BBGGRR) equivalents of the Office Theme hex colours
# 000000 / FFFFFF / 44546A / E7E6E6 / 5B9BD5 / ED7D31 / A5A5A5 / FFC000 /
# 4472C4 / 70AD47 / 0563C1 / 954F72
$officeThemeRgb = @(
    0,
    16777215,
    6968388,
    15132391,
    13998939,
    3243501,
    10855845,
    49407,
    12874308,
    4697456,
    12673797,
    7491477
)
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeThemeRgb[$i - 1]
}
